# Auto-generated edit script applying the Sagittarius_Profits market-price refresh
# (values for currentAveragePrice / LevePrice / LeveProfit columns per leve row).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 6.571429
$ws.Range("I11").Value = 6.571429
$ws.Range("K11").Value = 6.571429
$ws.Range("M11").Value = 133.428571

$ws.Range("H17").Value = 9634.923000000001
$ws.Range("I17").Value = 1573
$ws.Range("J17").Value = 11100.728
$ws.Range("K17").Value = 4719
$ws.Range("L17").Value = 33302.18399999999
$ws.Range("M17").Value = -4551
$ws.Range("N17").Value = -33638.18399999999

$ws.Range("H32").Value = 2319.1
$ws.Range("I32").Value = 800
$ws.Range("J32").Value = 2487.889
$ws.Range("K32").Value = 800
$ws.Range("L32").Value = 2487.889
$ws.Range("M32").Value = -474
$ws.Range("N32").Value = -3139.889

$ws.Range("H33").Value = 1196.3334
$ws.Range("I33").Value = 1587
$ws.Range("J33").Value = 415
$ws.Range("K33").Value = 1587
$ws.Range("L33").Value = 415
$ws.Range("M33").Value = -1358
$ws.Range("N33").Value = -873

$ws.Range("H41").Value = 1424.875
$ws.Range("I41").Value = 1057
$ws.Range("J41").Value = 4000
$ws.Range("K41").Value = 1057
$ws.Range("L41").Value = 4000
$ws.Range("M41").Value = -617
$ws.Range("N41").Value = -4880

$ws.Range("H53").Value = 383.9091
$ws.Range("I53").Value = 349.2
$ws.Range("J53").Value = 412.83334
$ws.Range("K53").Value = 349.2
$ws.Range("L53").Value = 412.83334
$ws.Range("M53").Value = 287.8
$ws.Range("N53").Value = -1686.83334

$ws.Range("H112").Value = 2543.5454
$ws.Range("J112").Value = 2469.2856
$ws.Range("L112").Value = 7407.8568
$ws.Range("N112").Value = -9623.856800000001

$ws.Range("H125").Value = 166684670
$ws.Range("I125").Value = 250014500
$ws.Range("J125").Value = 83354830
$ws.Range("K125").Value = 2250130500
$ws.Range("L125").Value = 750193470
$ws.Range("M125").Value = -2250128040
$ws.Range("N125").Value = -750198390

$ws.Range("H132").Value = 874.8333
$ws.Range("I132").Value = 969.8
$ws.Range("K132").Value = 2909.4
$ws.Range("M132").Value = -379.3999999999996

$ws.Range("H137").Value = 2006.6666
$ws.Range("I137").Value = 1616.6
$ws.Range("J137").Value = 2285.2856
$ws.Range("K137").Value = 4849.799999999999
$ws.Range("L137").Value = 6855.8568
$ws.Range("M137").Value = -2299.799999999999
$ws.Range("N137").Value = -11955.8568

$ws.Range("H141").Value = 1606.591
$ws.Range("I141").Value = 1587.8572
$ws.Range("J141").Value = 2000
$ws.Range("K141").Value = 4763.571599999999
$ws.Range("L141").Value = 6000
$ws.Range("M141").Value = 416.4284000000007
$ws.Range("N141").Value = -16360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2498.25
$ws.Range("I122").Value = 997
$ws.Range("J122").Value = 3999.5
$ws.Range("K122").Value = 2991
$ws.Range("L122").Value = 11998.5
$ws.Range("M122").Value = -541
$ws.Range("N122").Value = -16898.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1400
$ws.Range("I20").Value = 1400
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 1400
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()
$ws.Range("M20").Value = -1153

$ws.Range("H86").Value = 2211.6
$ws.Range("I86").Value = 2029
$ws.Range("J86").Value = 2333.3333
$ws.Range("K86").Value = 2029
$ws.Range("L86").Value = 2333.3333
$ws.Range("M86").Value = -906
$ws.Range("N86").Value = -4579.3333

$ws.Range("H89").Value = 2211.6
$ws.Range("I89").Value = 2029
$ws.Range("J89").Value = 2333.3333
$ws.Range("K89").Value = 10145
$ws.Range("L89").Value = 11666.6665
$ws.Range("M89").Value = -4529
$ws.Range("N89").Value = -22898.6665

$ws.Range("H99").Value = 1509.7273
$ws.Range("I99").Value = 1260.7
$ws.Range("J99").Value = 4000
$ws.Range("K99").Value = 1260.7
$ws.Range("L99").Value = 4000
$ws.Range("M99").Value = 237.3
$ws.Range("N99").Value = -6996

$ws.Range("H105").Value = 1632.85
$ws.Range("I105").Value = 1639.2941
$ws.Range("J105").Value = 1596.3334
$ws.Range("K105").Value = 1639.2941
$ws.Range("L105").Value = 1596.3334
$ws.Range("M105").Value = 107.7058999999999
$ws.Range("N105").Value = -5090.3334

$ws.Range("H115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("N115").ClearContents()
$ws.Range("L115").Value = 0

$ws.Range("H134").Value = 1444.1428
$ws.Range("I134").Value = 1551.5
$ws.Range("J134").Value = 800
$ws.Range("K134").Value = 4654.5
$ws.Range("L134").Value = 2400
$ws.Range("M134").Value = -2119.5
$ws.Range("N134").Value = -7470

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 5905.5
$ws.Range("J16").Value = 6196.3335
$ws.Range("L16").Value = 6196.3335
$ws.Range("N16").Value = -6770.3335

$ws.Range("H31").Value = 1589.7
$ws.Range("I31").Value = 1566.3334
$ws.Range("K31").Value = 1566.3334
$ws.Range("M31").Value = -1271.3334

$ws.Range("H34").Value = 1589.7
$ws.Range("I34").Value = 1566.3334
$ws.Range("K34").Value = 1566.3334
$ws.Range("M34").Value = -1364.3334

$ws.Range("H58").Value = 1229.96
$ws.Range("I58").Value = 1217.2273
$ws.Range("J58").Value = 1323.3334
$ws.Range("K58").Value = 1217.2273
$ws.Range("L58").Value = 1323.3334
$ws.Range("M58").Value = -1014.2273
$ws.Range("N58").Value = -1729.3334

$ws.Range("H113").Value = 5905.5
$ws.Range("J113").Value = 6196.3335
$ws.Range("L113").Value = 6196.3335
$ws.Range("N113").Value = -10536.3335

$ws.Range("H122").Value = 1974.8
$ws.Range("I122").Value = 2001.5714
$ws.Range("K122").Value = 6004.7142
$ws.Range("M122").Value = -3554.7142

$ws.Range("H132").Value = 1947.5625
$ws.Range("I132").Value = 1954.3572
$ws.Range("J132").Value = 1900
$ws.Range("K132").Value = 5863.071599999999
$ws.Range("L132").Value = 5700
$ws.Range("M132").Value = -3333.071599999999
$ws.Range("N132").Value = -10760

$ws.Range("H134").Value = 1374
$ws.Range("I134").Value = 1374
$ws.Range("K134").Value = 4122
$ws.Range("M134").Value = -1587

$ws.Range("H136").Value = 1229.96
$ws.Range("I136").Value = 1217.2273
$ws.Range("J136").Value = 1323.3334
$ws.Range("K136").Value = 3651.6819
$ws.Range("L136").Value = 3970.0002
$ws.Range("M136").Value = -1101.6819
$ws.Range("N136").Value = -9070.0002

$ws.Range("H141").Value = 195714.28
$ws.Range("J141").Value = 195714.28
$ws.Range("L141").Value = 195714.28
$ws.Range("N141").Value = -206074.28

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 5483.2856
$ws.Range("J39").Value = 5996.8
$ws.Range("L39").Value = 17990.4
$ws.Range("N39").Value = -18578.4

$ws.Range("H40").Value = 96.75
$ws.Range("I40").Value = 37.2
$ws.Range("K40").Value = 148.8
$ws.Range("M40").Value = -79.80000000000001

$ws.Range("H55").Value = 1495
$ws.Range("I55").Value = 1495
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 4485
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()
$ws.Range("M55").Value = -4308

$ws.Range("H131").Value = 419942.78
$ws.Range("J131").Value = 559328.9
$ws.Range("L131").Value = 1677986.7
$ws.Range("N131").Value = -1688066.7

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6698.7856
$ws.Range("I70").Value = 6334.5557
$ws.Range("J70").Value = 7354.4
$ws.Range("K70").Value = 6334.5557
$ws.Range("L70").Value = 7354.4
$ws.Range("M70").Value = -6064.5557
$ws.Range("N70").Value = -7894.4

$ws.Range("H73").Value = 6698.7856
$ws.Range("I73").Value = 6334.5557
$ws.Range("J73").Value = 7354.4
$ws.Range("K73").Value = 6334.5557
$ws.Range("L73").Value = 7354.4
$ws.Range("M73").Value = -5398.5557
$ws.Range("N73").Value = -9226.4

$ws.Range("H102").Value = 1269.2727
$ws.Range("I102").Value = 1269.2727
$ws.Range("K102").Value = 1269.2727
$ws.Range("M102").Value = 352.7273

$ws.Range("H132").Value = 1335.125
$ws.Range("I132").Value = 1335.125
$ws.Range("K132").Value = 4005.375
$ws.Range("M132").Value = -1475.375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1893.8823
$ws.Range("I22").Value = 1824.875
$ws.Range("J22").Value = 2998
$ws.Range("K22").Value = 1824.875
$ws.Range("L22").Value = 2998
$ws.Range("M22").Value = -1529.875
$ws.Range("N22").Value = -3588

$ws.Range("H27").Value = 1893.8823
$ws.Range("I27").Value = 1824.875
$ws.Range("J27").Value = 2998
$ws.Range("K27").Value = 1824.875
$ws.Range("L27").Value = 2998
$ws.Range("M27").Value = -1717.875
$ws.Range("N27").Value = -3212

$ws.Range("H40").Value = 2219.3
$ws.Range("I40").Value = 1846.2354
$ws.Range("K40").Value = 1846.2354
$ws.Range("M40").Value = -1710.2354

$ws.Range("H55").Value = 215
$ws.Range("I55").Value = 63.333332
$ws.Range("J55").Value = 366.66666
$ws.Range("K55").Value = 63.333332
$ws.Range("L55").Value = 366.66666
$ws.Range("M55").Value = 109.666668
$ws.Range("N55").Value = -712.66666

$ws.Range("H93").Value = 2433.3333
$ws.Range("J93").Value = 2360
$ws.Range("L93").Value = 2360
$ws.Range("N93").Value = -4856

$ws.Range("H136").Value = 3083.913
$ws.Range("I136").Value = 2646.35
$ws.Range("J136").Value = 6001
$ws.Range("K136").Value = 7939.049999999999
$ws.Range("L136").Value = 18003
$ws.Range("M136").Value = -5389.049999999999
$ws.Range("N136").Value = -23103

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1051.5
$ws.Range("I96").Value = 803
$ws.Range("J96").Value = 1300
$ws.Range("K96").Value = 803
$ws.Range("L96").Value = 1300
$ws.Range("M96").Value = 570
$ws.Range("N96").Value = -4046

$ws.Range("H100").Value = 5883219.5
$ws.Range("I100").Value = 6250902
$ws.Range("J100").Value = 299
$ws.Range("K100").Value = 12501804
$ws.Range("L100").Value = 598
$ws.Range("M100").Value = -12501263
$ws.Range("N100").Value = -1680

$ws.Range("H132").Value = 2790.65
$ws.Range("I132").Value = 3007
$ws.Range("J132").Value = 1925.25
$ws.Range("K132").Value = 9021
$ws.Range("L132").Value = 5775.75
$ws.Range("M132").Value = -6491
$ws.Range("N132").Value = -10835.75
